$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.686.26'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.596.67'
$ws.Range("E3").Value = '  -1.15%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '211.44'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.31%  '
$ws.Range("E6").Value = '  +0.17%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -0.78%  '
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.72'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.30%  '
$ws.Range("E11").Value = '  +0.28%  '
$ws.Range("D12").Value = '1.821.27'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("D13").Value = '1.602.46'
$ws.Range("E13").Value = '  -0.73%  '
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("E15").Value = '  -2.01%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.08'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.99%  '
$ws.Range("D17").Value = '26.690.54'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").Value = '0.0₃0728'
$ws.Range("E18").Value = '  -0.35%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '210.12'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  -0.37%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.26'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.93%  '
$ws.Range("E23").Value = '  -1.42%  '
$ws.Range("E24").Value = '  -0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.16%  '
$ws.Range("E26").Value = '  +0.08%  '
$ws.Range("E27").Value = '  -3.88%  '
$ws.Range("E28").Value = '  +1.94%  '
$ws.Range("E29").Value = '  -0.59%  '
$ws.Range("E30").Value = '  +0.14%  '
$ws.Range("E31").Value = '  -0.69%  '
$ws.Range("E32").Value = '  -1.39%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.669'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.51%  '
$ws.Range("E34").Value = '  -1.89%  '
$ws.Range("D35").Value = '1.295.16'
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("E37").Value = '  -4.27%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0171'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.54%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.843'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.65%  '
$ws.Range("E41").Value = '  -0.49%  '
$ws.Range("E42").Value = '  +1.42%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '63.58'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.22%  '
$ws.Range("D45").Value = '1.733.04'
$ws.Range("E45").Value = '  -0.99%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '90.16'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.85%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.879'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +8.36%  '
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("E49").Value = '  +1.27%  '
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.50'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.42%  '
